# edit.ps1 - apply the CV 2.0 Rafal Mucha changes described in the commit:
#   "Added Ionic cert, seperated CSS & JS, fixed btns, added new effects"
#
# Concrete, controllable content edits:
#   1. Skills list: "TypeScript, SASS" -> "TypeScript, SASS, Figma"
#   2. Tech badge table: "Node" + "JS" (two runs) -> single "NodeJS" run
#   3. Tech badge table: "Express" + "JS" (two runs) -> single "ExpressJS" run
#   4. Tech badge table: "Ruby" -> "Ionic" (new certification badge)

$d = $word.ActiveDocument

# wdReplaceOne / wdFindWrapStop constants used below via Find.Execute's
# numeric Replace parameter (1 = wdReplaceOne) and Wrap parameter (0 = wdFindStop).
$wdReplaceOne = 1
$wdFindStop = 0

# --- 1. Add ", Figma" to the "TypeScript, SASS" skill line -----------------
$rng = $d.Content
$found = $rng.Find.Execute("TypeScript, SASS", $true, $true, $false, $false, $false, $true, $wdFindStop, $false, "", 0)
if ($found) {
    $rng.Collapse(0)   # wdCollapseEnd - move to just after the matched text
    $rng.InsertAfter(", Figma")
}

# --- 2. Merge "Node" + "JS" runs into "NodeJS" ------------------------------
$d.Content.Find.Execute("Node" + "JS", $true, $true, $false, $false, $false, $true, $wdFindStop, $false, "NodeJS", 2) | Out-Null

# --- 3. Merge "Express" + "JS" runs into "ExpressJS" ------------------------
$d.Content.Find.Execute("Express" + "JS", $true, $true, $false, $false, $false, $true, $wdFindStop, $false, "ExpressJS", 2) | Out-Null

# --- 4. Replace the "Ruby" badge with the new "Ionic" certification --------
$d.Content.Find.Execute("Ruby", $true, $true, $false, $false, $false, $true, $wdFindStop, $false, "Ionic", 2) | Out-Null
